$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("USERS")

$ws.Range("B9").Value = "displayName"
$ws.Range("C9").Value = "string"
$ws.Range("D9").Value = "닉네임"

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:F9"))

$ws.Range("D9").Select()
